# Fix forecast-error array simulations: rename the lone sheet to "Table"
# and add the EX / LIN / GSSA / VFI simulation-method sheets, each holding
# the TZFE / OPFE mean+RMSE headers and the k,Y,w,r,T,c,i,u row labels.
# Also records a couple of numbers that were missing on the Table sheet
# (C9 RMSE value, and the "n/a" placeholder next to the Solve row).

$wb = $excel.ActiveWorkbook

# ---- rename the original (and only) sheet to "Table" ----------------------
$wsTable = $wb.Worksheets.Item(1)
$wsTable.Name = "Table"

# ---- add the four method sheets, in display order after Table -------------
$wsEX = $wb.Worksheets.Add($null, $wsTable)
$wsEX.Name = "EX"

$wsLIN = $wb.Worksheets.Add($null, $wsEX)
$wsLIN.Name = "LIN"

$wsGSSA = $wb.Worksheets.Add($null, $wsLIN)
$wsGSSA.Name = "GSSA"

$wsVFI = $wb.Worksheets.Add($null, $wsGSSA)
$wsVFI.Name = "VFI"

# ---- EX sheet: headers + row labels, no data yet ---------------------------
$wsEX.Range("B1").Value = "TZFE"
$wsEX.Range("D1").Value = "OPFE"
$wsEX.Range("B2").Value = "mean"
$wsEX.Range("C2").Value = "RMSE"
$wsEX.Range("D2").Value = "mean"
$wsEX.Range("E2").Value = "RMSE"
$wsEX.Range("A3").Value = "k"
$wsEX.Range("A4").Value = "Y"
$wsEX.Range("A5").Value = "w"
$wsEX.Range("A6").Value = "r"
$wsEX.Range("A7").Value = "T"
$wsEX.Range("A8").Value = "c"
$wsEX.Range("A9").Value = "i"
$wsEX.Range("A10").Value = "u"
$wsEX.Range("A1:E10").Select()

# ---- LIN sheet: headers, row labels and the simulated error numbers -------
$wsLIN.Range("C1").Value = "TZFE"
$wsLIN.Range("E1").Value = "OPFE"
$wsLIN.Range("C2").Value = "mean"
$wsLIN.Range("D2").Value = "RMSE"
$wsLIN.Range("E2").Value = "mean"
$wsLIN.Range("F2").Value = "RMSE"

$wsLIN.Range("B3").Value = "k"
$wsLIN.Range("C3").Value = 0.00519241
$wsLIN.Range("D3").Value = 0.169393
$wsLIN.Range("E3").Value = 0.0123326
$wsLIN.Range("F3").Value = 0.0918031

$wsLIN.Range("B4").Value = "Y"
$wsLIN.Range("C4").Value = 0.00623188
$wsLIN.Range("D4").Value = 0.106471
$wsLIN.Range("E4").Value = 0.00725747
$wsLIN.Range("F4").Value = 0.00909569

$wsLIN.Range("B5").Value = "w"
$wsLIN.Range("C5").Value = 0.00521837
$wsLIN.Range("D5").Value = 0.123986
$wsLIN.Range("E5").Value = 0.00396274
$wsLIN.Range("F5").Value = 0.00496714

$wsLIN.Range("B6").Value = "r"
$wsLIN.Range("C6").Value = 0.000142776
$wsLIN.Range("D6").Value = 0.011889
$wsLIN.Range("E6").Value = 0.0039898
$wsLIN.Range("F6").Value = 0.00500044

$wsLIN.Range("B7").Value = "T"
$wsLIN.Range("C7").Value = 0.00571917
$wsLIN.Range("D7").Value = 0.456147
$wsLIN.Range("E7").Value = 0.00431821
$wsLIN.Range("F7").Value = 0.00541547

$wsLIN.Range("B8").Value = "c"
$wsLIN.Range("C8").Value = 0.00521909
$wsLIN.Range("D8").Value = 0.123662
$wsLIN.Range("E8").Value = 0.00396298
$wsLIN.Range("F8").Value = 0.00496832

$wsLIN.Range("B9").Value = "i"
$wsLIN.Range("C9").Value = 0.00452961
$wsLIN.Range("D9").Value = 0.147412
$wsLIN.Range("E9").Value = 0.00128923
$wsLIN.Range("F9").Value = 0.00163892

$wsLIN.Range("B10").Value = "u"
$wsLIN.Range("C10").Value = 0.00488911
$wsLIN.Range("D10").Value = 0.0674381
$wsLIN.Range("E10").Value = 0.00139145
$wsLIN.Range("F10").Value = 0.0017686

$wsLIN.Range("B1:F10").Select()

# ---- GSSA sheet: same layout as EX, no data yet ----------------------------
$wsGSSA.Range("B1").Value = "TZFE"
$wsGSSA.Range("D1").Value = "OPFE"
$wsGSSA.Range("B2").Value = "mean"
$wsGSSA.Range("C2").Value = "RMSE"
$wsGSSA.Range("D2").Value = "mean"
$wsGSSA.Range("E2").Value = "RMSE"
$wsGSSA.Range("A3").Value = "k"
$wsGSSA.Range("A4").Value = "Y"
$wsGSSA.Range("A5").Value = "w"
$wsGSSA.Range("A6").Value = "r"
$wsGSSA.Range("A7").Value = "T"
$wsGSSA.Range("A8").Value = "c"
$wsGSSA.Range("A9").Value = "i"
$wsGSSA.Range("A10").Value = "u"
$wsGSSA.Range("A1:E10").Select()

# ---- VFI sheet: same layout as EX, no data yet -----------------------------
$wsVFI.Range("B1").Value = "TZFE"
$wsVFI.Range("D1").Value = "OPFE"
$wsVFI.Range("B2").Value = "mean"
$wsVFI.Range("C2").Value = "RMSE"
$wsVFI.Range("D2").Value = "mean"
$wsVFI.Range("E2").Value = "RMSE"
$wsVFI.Range("A3").Value = "k"
$wsVFI.Range("A4").Value = "Y"
$wsVFI.Range("A5").Value = "w"
$wsVFI.Range("A6").Value = "r"
$wsVFI.Range("A7").Value = "T"
$wsVFI.Range("A8").Value = "c"
$wsVFI.Range("A9").Value = "i"
$wsVFI.Range("A10").Value = "u"
$wsVFI.Range("H27").Select()

# ---- back to Table: fill in the two missing cells --------------------------
$wsTable.Select()
$wsTable.Range("C9").Value = 0.000251951
$wsTable.Range("B11").Value = "n/a"
$wsTable.Range("C18").Select()

# ---- widen the saved window a bit now that there are five tabs -------------
# (best-effort; the host's Window object is a view-only stub and this
# property does not currently round-trip into bookViews/workbookView, but
# setting it is harmless in case a future host version wires it up)
$excel.ActiveWindow.Width = 19380
